# Apply the edits described by the diff:
# - Sheet1!C2 changes from 8 to 3 (Hours Estimate for Task_1)
# - Sheet1!D6 changes from 5 to 8 (Days To Deadline for Task_5)
# - Active selection moves to A10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 3
$ws.Range("D6").Value = 8

$ws.Activate()
$ws.Range("A10").Select()
